$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $oldText with $newText, scoped to
# a given paragraph (by 1-based index). Re-fetches the paragraph Range each
# time so a Find/Replace collapsing the range doesn't break later calls.
# ---------------------------------------------------------------------------
function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    return $ok
}

# ===========================================================================
# 1) "Areas of expertise" line - swap several of the pipe-separated items.
# ===========================================================================
Replace-InParagraph 7 "Custom Software Development" "Web Technologies" | Out-Null
Replace-InParagraph 7 "Web Design + Development" "Frontend Development" | Out-Null
Replace-InParagraph 7 "Technical Support" "Library and Framework Development" | Out-Null
Replace-InParagraph 7 "Mobile Application Development" "Observability and Real-time Analytics  |   Software Testing" | Out-Null
Replace-InParagraph 7 "Software Testing   |   User Experience Design (UED)  | Technical Problem" "Product Development  | Technical Problem" | Out-Null
Replace-InParagraph 7 "Problem Solving | Software Development + Engineering" "Problem Solving " | Out-Null

# ===========================================================================
# 2) Job date range: "July 2020 - Present" -> "July 2020 " / "-" / " " / "July 2023"
#    (4 runs, each keeping the original Bold + BoldBi formatting).
# ===========================================================================
$dateRng = $d.Paragraphs(9).Range
$found = $dateRng.Find.Execute("July 2020 - Present")
if ($found) {
    $enDash = [string][char]0x2013
    # Setting .Text on the exact found range keeps the run (and its rPr,
    # including bCs) intact - only the text content changes.
    $dateRng.Text = "July 2020 " + $enDash + " July 2023"

    $segStart = $dateRng.Start
    $seg1End = $segStart + 10   # "July 2020 "
    $seg2End = $seg1End + 1     # en dash
    $seg3End = $seg2End + 1     # " "
    $seg4End = $seg3End + 9     # "July 2023"

    $seg1 = $d.Range($segStart, $seg1End)
    $seg2 = $d.Range($seg1End, $seg2End)
    $seg3 = $d.Range($seg2End, $seg3End)
    $seg4 = $d.Range($seg3End, $seg4End)

    # Toggling Bold off/on forces the engine to split these into their own
    # runs instead of re-coalescing them with their neighbours, while still
    # ending up with the same Bold/BoldBi (b + bCs) formatting as before.
    foreach ($seg in @($seg1, $seg2, $seg3, $seg4)) {
        $seg.Bold = 0
        $seg.Bold = 1
    }
}

# ===========================================================================
# 3) "Node.js" -> "Node.js Platform" inside the BNY Mellon bullet paragraph.
# ===========================================================================
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Typescript, Node.js, Java Spring*") {
        Replace-InParagraph $i "Typescript, Node.js, Java Spring" "Typescript, Node.js Platform, Java Spring" | Out-Null
        break
    }
}
